$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Update the "Scenario Indices"/block cell references in column D (rows 5-11)
# from referencing row 26 to referencing row 28 (scenario 25 -> scenario 27 update).
$ws.Range("D5").Value = "A28"
$ws.Range("D6").Value = "B28"
$ws.Range("D7").Value = "C28"
$ws.Range("D8").Value = "G28"
$ws.Range("D9").Value = "H28"
$ws.Range("D10").Value = "I28"
$ws.Range("D11").Value = "J28"

# Move the active selection to D5, matching the saved cursor position.
$ws.Range("D5").Select()
